$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 515374.88
$ws.Range("J43").Value = 588857.1
$ws.Range("L43").Value = 588857.1
$ws.Range("N43").Value = -588995.1

$ws.Range("H76").Value = 7162.5835
$ws.Range("I76").Value = 6177.636
$ws.Range("J76").Value = 17997
$ws.Range("K76").Value = 6177.636
$ws.Range("L76").Value = 17997
$ws.Range("M76").Value = -5862.636
$ws.Range("N76").Value = -18627

$ws.Range("H79").Value = 7162.5835
$ws.Range("I79").Value = 6177.636
$ws.Range("J79").Value = 17997
$ws.Range("K79").Value = 6177.636
$ws.Range("L79").Value = 17997
$ws.Range("M79").Value = -5085.636
$ws.Range("N79").Value = -20181

$ws.Range("H113").Value = 29426362
$ws.Range("J113").Value = 41686308
$ws.Range("L113").Value = 41686308
$ws.Range("N113").Value = -41692816

$ws.Range("H121").Value = 1500
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994

$ws.Range("H135").Value = 833960.7
$ws.Range("I135").Value = 909590.2
$ws.Range("K135").Value = 8186311.8
$ws.Range("M135").Value = -8183776.8

$ws.Range("H138").Value = 1355665.5
$ws.Range("I138").Value = 2574.3809
$ws.Range("J138").Value = 1891796
$ws.Range("K138").Value = 7723.1427
$ws.Range("L138").Value = 5675388
$ws.Range("M138").Value = -2583.1427
$ws.Range("N138").Value = -5685668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 3694.5454
$ws.Range("I16").Value = 352.5
$ws.Range("J16").Value = 12606.667
$ws.Range("K16").Value = 352.5
$ws.Range("L16").Value = 12606.667
$ws.Range("M16").Value = -65.5
$ws.Range("N16").Value = -13180.667

$ws.Range("H32").Value = 3781201.2
$ws.Range("I32").Value = 4006073.2
$ws.Range("K32").Value = 4006073.2
$ws.Range("M32").Value = -4005786.2

$ws.Range("H61").Value = 27028770
$ws.Range("I61").Value = 1202.25
$ws.Range("J61").Value = 200005200
$ws.Range("K61").Value = 1202.25
$ws.Range("L61").Value = 200005200
$ws.Range("M61").Value = -990.25
$ws.Range("N61").Value = -200005624

$ws.Range("H74").Value = 25086.068
$ws.Range("I74").Value = 37264.395
$ws.Range("K74").Value = 37264.395
$ws.Range("M74").Value = -36390.395

$ws.Range("H77").Value = 25086.068
$ws.Range("I77").Value = 37264.395
$ws.Range("K77").Value = 186321.975
$ws.Range("M77").Value = -181953.975

$ws.Range("H102").Value = 3815.4707
$ws.Range("I102").Value = 2905.3333
$ws.Range("K102").Value = 2905.3333
$ws.Range("M102").Value = -1283.3333

$ws.Range("H103").Value = 53737.5
$ws.Range("J103").Value = 53737.5
$ws.Range("L103").Value = 53737.5
$ws.Range("N103").Value = -56081.5

$ws.Range("H122").Value = 2839.878
$ws.Range("I122").Value = 2496.6365
$ws.Range("J122").Value = 4255.75
$ws.Range("K122").Value = 7489.9095
$ws.Range("L122").Value = 12767.25
$ws.Range("M122").Value = -5039.9095
$ws.Range("N122").Value = -17667.25

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H136").Value = 27028770
$ws.Range("I136").Value = 1202.25
$ws.Range("J136").Value = 200005200
$ws.Range("K136").Value = 3606.75
$ws.Range("L136").Value = 600015600
$ws.Range("M136").Value = -1056.75
$ws.Range("N136").Value = -600020700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6944712
$ws.Range("I22").Value = 9259509
$ws.Range("J22").Value = 322.5
$ws.Range("K22").Value = 9259509
$ws.Range("L22").Value = 322.5
$ws.Range("M22").Value = -9259336
$ws.Range("N22").Value = -668.5

$ws.Range("H107").Value = 51137380
$ws.Range("I107").Value = 59211548
$ws.Range("J107").Value = 965.3333
$ws.Range("K107").Value = 59211548
$ws.Range("L107").Value = 965.3333
$ws.Range("M107").Value = -59209628
$ws.Range("N107").Value = -4805.3333

$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 49999
$ws.Range("K129").Value = 49999
$ws.Range("M129").Value = -44999

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws.Range("H134").Value = 7870.6
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 568
$ws.Range("I14").Value = 568
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 568
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -398
$ws.Range("N14").Value = $null

$ws.Range("H31").Value = 4904.959
$ws.Range("I31").Value = 1754.45
$ws.Range("K31").Value = 1754.45
$ws.Range("M31").Value = -1459.45

$ws.Range("H34").Value = 4904.959
$ws.Range("I34").Value = 1754.45
$ws.Range("K34").Value = 1754.45
$ws.Range("M34").Value = -1552.45

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H132").Value = 3920.9375
$ws.Range("I132").Value = 2693.6875
$ws.Range("K132").Value = 8081.0625
$ws.Range("M132").Value = -5551.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = $null

$ws.Range("H131").Value = 1727.1875
$ws.Range("J131").Value = 1803.6154
$ws.Range("L131").Value = 5410.8462
$ws.Range("N131").Value = -15490.8462

$ws.Range("H132").Value = 5573.343
$ws.Range("I132").Value = 2624.5
$ws.Range("J132").Value = 7111.8696
$ws.Range("K132").Value = 23620.5
$ws.Range("L132").Value = 64006.8264
$ws.Range("M132").Value = -21090.5
$ws.Range("N132").Value = -69066.82639999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 471280.06
$ws.Range("I107").Value = 1000251.5
$ws.Range("K107").Value = 1000251.5
$ws.Range("M107").Value = -998331.5

$ws.Range("H122").Value = 4542460
$ws.Range("I122").Value = 6604395
$ws.Range("K122").Value = 19813185
$ws.Range("M122").Value = -19810735

$ws.Range("H129").Value = 59480
$ws.Range("J129").Value = 59480
$ws.Range("L129").Value = 59480
$ws.Range("N129").Value = -69480

$ws.Range("H132").Value = 3267.5386
$ws.Range("I132").Value = 2318.4666
$ws.Range("J132").Value = 4561.727
$ws.Range("K132").Value = 6955.399800000001
$ws.Range("L132").Value = 13685.181
$ws.Range("M132").Value = -4425.399800000001
$ws.Range("N132").Value = -18745.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 40000
$ws.Range("J25").Value = 40000
$ws.Range("L25").Value = 40000
$ws.Range("N25").Value = -40460

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null

$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = $null

$ws.Range("H46").Value = 1612.2727
$ws.Range("I46").Value = 1369.8948
$ws.Range("K46").Value = 1369.8948
$ws.Range("M46").Value = -1181.8948

$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = $null

$ws.Range("H116").Value = 56067
$ws.Range("J116").Value = 56067
$ws.Range("L116").Value = 56067
$ws.Range("N116").Value = -65245

$ws.Range("H122").Value = 5408.095
$ws.Range("I122").Value = 4956.9
$ws.Range("J122").Value = 5818.273
$ws.Range("K122").Value = 14870.7
$ws.Range("L122").Value = 17454.819
$ws.Range("M122").Value = -12420.7
$ws.Range("N122").Value = -22354.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 8999.666999999999
$ws.Range("J20").Value = 8999.666999999999
$ws.Range("L20").Value = 8999.666999999999
$ws.Range("N20").Value = -9479.666999999999

$ws.Range("H107").Value = 22223556
$ws.Range("I107").Value = 406
$ws.Range("K107").Value = 1218
$ws.Range("M107").Value = 702

$ws.Range("H132").Value = 3943.9546
$ws.Range("I132").Value = 3708.6287
$ws.Range("K132").Value = 11125.8861
$ws.Range("M132").Value = -8595.8861
